$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$b2Text = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Akita']
MATCH (c)<--(diag:diagnosis)
 MATCH (co:cohort)<-[*]-(c)
 OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID`,
        coalesce(s.clinical_study_designation, '') AS `Study Code`,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
        coalesce(demo.patient_age_at_enrollment, '') AS Age,
        coalesce(demo.sex, '') AS Sex,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
        coalesce(co.cohort_description, '') AS `Cohort`
'@

$c2Text = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Akita']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN  count(DISTINCT(s.clinical_study_designation)) as Studies, count(DISTINCT(c.case_id)) as Cases ,count(DISTINCT(samp)) as Samples ,count(DISTINCT(f)) as `Case Files` 
'@

$ws.Range("C2").Value = $c2Text
$ws.Range("B2").Value = $b2Text

$ws.Range("C2").Select()
